$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at position 178 (shifts existing rows 178-256 down to 183-261)
$ws.Rows("178:182").Insert()

# Row 178
$ws.Range("A178").Value = 'Paul Luo Li and Randy Nakagawa and Rob Montroy'
$ws.Range("B178").Value = '@inproceedings{luo:2007,
booktitle = {International Symposium on Empirical Software Engineering and Measurement (ESEM)}, 
author = {Paul Luo Li and Randy Nakagawa and Rob Montroy}, 
title = {{Estimating the Quality of Widely Used Software Products Using Software Reliability Growth Modeling: Case Study of an IBM Federated Database Project}}, 
year = {2007} 
}'
$ws.Range("C178").Value = 'conf/esem/LiNM07'
$ws.Range("D178").Value = 'Estimating the Quality of Widely Used Software Products Using Software Reliability Growth Modeling: Case Study of an IBM Federated Database Project'
$ws.Range("E178").Value = 'db/conf/esem/esem2007.html#LiNM07'
$ws.Range("F178").Value = 'conf/esem/'
$ws.Range("H178").Value = 'conf/esem/2007'
$ws.Range("G178").Value = 2007
# Reset auto-calculated row height (from multi-line bibtex text) back to sheet default
$ws.Rows("178").AutoFit()

# Row 179
$ws.Range("A179").Value = 'Donald W. McCormick II and William B. Frakes and Reghu Anguswamy'
$ws.Range("B179").Value = '@inproceedings{w.:2012,
booktitle = {International Symposium on Empirical Software Engineering and Measurement (ESEM)}, 
author = {Donald W. McCormick II and William B. Frakes and Reghu Anguswamy}, 
title = {{A comparison of database fault detection capabilities using mutation testing}}, 
year = {2012} 
}'
$ws.Range("C179").Value = 'conf/esem/McCormickFA12'
$ws.Range("D179").Value = 'A comparison of database fault detection capabilities using mutation testing'
$ws.Range("E179").Value = 'db/conf/esem/esem2012.html#McCormickFA12'
$ws.Range("F179").Value = 'conf/esem/'
$ws.Range("H179").Value = 'conf/esem/2012'
$ws.Range("G179").Value = 2012
# Reset auto-calculated row height (from multi-line bibtex text) back to sheet default
$ws.Rows("179").AutoFit()

# Row 180
$ws.Range("A180").Value = 'Rudolf Ramler and Klaus Wolfmaier'
$ws.Range("B180").Value = '@inproceedings{ramler:2008,
booktitle = {International Symposium on Empirical Software Engineering and Measurement (ESEM)}, 
author = {Rudolf Ramler and Klaus Wolfmaier}, 
title = {{Issues and effort in integrating data from heterogeneous software repositories and corporate databases}}, 
year = {2008} 
}'
$ws.Range("C180").Value = 'conf/esem/RamlerW08'
$ws.Range("D180").Value = 'Issues and effort in integrating data from heterogeneous software repositories and corporate databases'
$ws.Range("E180").Value = 'db/conf/esem/esem2008.html#RamlerW08'
$ws.Range("F180").Value = 'conf/esem/'
$ws.Range("H180").Value = 'conf/esem/2008'
$ws.Range("G180").Value = 2008
# Reset auto-calculated row height (from multi-line bibtex text) back to sheet default
$ws.Rows("180").AutoFit()

# Row 181
$ws.Range("A181").Value = 'Michael Wedel and Uwe Jensen and Peter Göhner'
$ws.Range("B181").Value = '@inproceedings{wedel:2008,
booktitle = {International Symposium on Empirical Software Engineering and Measurement (ESEM)}, 
author = {Michael Wedel and Uwe Jensen and Peter Göhner}, 
title = {{Mining software code repositories and bug databases using survival analysis models}}, 
year = {2008} 
}'
$ws.Range("C181").Value = 'conf/esem/WedelJG08'
$ws.Range("D181").Value = 'Mining software code repositories and bug databases using survival analysis models'
$ws.Range("E181").Value = 'db/conf/esem/esem2008.html#WedelJG08'
$ws.Range("F181").Value = 'conf/esem/'
$ws.Range("H181").Value = 'conf/esem/2008'
$ws.Range("G181").Value = 2008
# Reset auto-calculated row height (from multi-line bibtex text) back to sheet default
$ws.Rows("181").AutoFit()

# Row 182
$ws.Range("A182").Value = 'Samireh Jalali and Claes Wohlin'
$ws.Range("B182").Value = '@inproceedings{jalali:2012,
booktitle = {International Symposium on Empirical Software Engineering and Measurement (ESEM)}, 
author = {Samireh Jalali and Claes Wohlin}, 
title = {{Systematic literature studies: database searches vs. backward snowballing}}, 
year = {2012} 
}'
$ws.Range("C182").Value = 'conf/esem/JalaliW12'
$ws.Range("D182").Value = 'Systematic literature studies: database searches vs. backward snowballing'
$ws.Range("E182").Value = 'db/conf/esem/esem2012.html#JalaliW12'
$ws.Range("F182").Value = 'conf/esem/'
$ws.Range("H182").Value = 'conf/esem/2012'
$ws.Range("G182").Value = 2012
# Reset auto-calculated row height (from multi-line bibtex text) back to sheet default
$ws.Rows("182").AutoFit()
